$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.151.67'
$ws.Cells.Item(2, 5).Value = '  -2.21%  '
$ws.Cells.Item(3, 4).Value = '1.574.71'
$ws.Cells.Item(3, 5).Value = '  -1.76%  '
$ws.Cells.Item(4, 5).Value = '  -0.47%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '208.73'
$ws.Cells.Item(5, 5).Value = '  -1.45%  '
$ws.Cells.Item(6, 5).Value = '  -3.03%  '
$ws.Cells.Item(7, 5).Value = '  -0.43%  '
$ws.Cells.Item(8, 5).Value = '  -1.69%  '
$ws.Cells.Item(9, 5).Value = '  -1.31%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '19.57'
$ws.Cells.Item(10, 5).Value = '  -0.73%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0844'
$ws.Cells.Item(11, 5).Value = '  -0.37%  '
$ws.Cells.Item(12, 4).Value = '1.796.14'
$ws.Cells.Item(12, 5).Value = '  -1.76%  '
$ws.Cells.Item(13, 2).Value = 'Polkadot'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.06'
$ws.Cells.Item(13, 5).Value = '  -0.45%  '
$ws.Cells.Item(14, 2).Value = 'WrappedEther'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(14, 4).Value = '1.565.71'
$ws.Cells.Item(14, 5).Value = '  -2.09%  '
$ws.Cells.Item(15, 5).Value = '  -2.12%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '64.34'
$ws.Cells.Item(16, 5).Value = '  -1.10%  '
$ws.Cells.Item(17, 4).Value = '26.138.33'
$ws.Cells.Item(17, 5).Value = '  -2.15%  '
$ws.Cells.Item(18, 5).Value = '  -2.40%  '
$ws.Cells.Item(19, 5).Value = '  +1.90%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '207.64'
$ws.Cells.Item(20, 5).Value = '  -1.08%  '
$ws.Cells.Item(21, 5).Value = '  -0.40%  '
$ws.Cells.Item(22, 5).Value = '  -1.19%  '
$ws.Cells.Item(23, 5).Value = '  -2.73%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '8.82'
$ws.Cells.Item(24, 5).Value = '  -2.94%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '143.68'
$ws.Cells.Item(25, 5).Value = '  +0.00%  '
$ws.Cells.Item(26, 5).Value = '  -0.51%  '
$ws.Cells.Item(27, 5).Value = '  -1.68%  '
$ws.Cells.Item(29, 5).Value = '  -1.13%  '
$ws.Cells.Item(30, 5).Value = '  -0.56%  '
$ws.Cells.Item(31, 5).Value = '  -1.54%  '
$ws.Cells.Item(32, 5).Value = '  -2.04%  '
$ws.Cells.Item(33, 5).Value = '  +0.26%  '
$ws.Cells.Item(34, 4).Value = '1.276.78'
$ws.Cells.Item(34, 5).Value = '  -1.05%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.619'
$ws.Cells.Item(35, 5).Value = '  +4.43%  '
$ws.Cells.Item(36, 5).Value = '  -1.66%  '
$ws.Cells.Item(37, 5).Value = '  -0.97%  '
$ws.Cells.Item(38, 2).Value = 'WEMIXToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '1.11'
$ws.Cells.Item(38, 5).Value = '  -10.35%  '
$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.0166'
$ws.Cells.Item(39, 5).Value = '  -2.59%  '
$ws.Cells.Item(40, 5).Value = '  -2.31%  '
$ws.Cells.Item(41, 5).Value = '  -0.41%  '
$ws.Cells.Item(42, 5).Value = '  +2.17%  '
$ws.Cells.Item(43, 5).Value = '  -2.63%  '
$ws.Cells.Item(44, 5).Value = '  -2.12%  '
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '62.27'
$ws.Cells.Item(45, 5).Value = '  -1.07%  '
$ws.Cells.Item(46, 4).Value = '1.709.06'
$ws.Cells.Item(46, 5).Value = '  -1.78%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '88.64'
$ws.Cells.Item(47, 5).Value = '  -2.05%  '
$ws.Cells.Item(48, 5).Value = '  -3.13%  '
$ws.Cells.Item(49, 5).Value = '  +0.55%  '
$ws.Cells.Item(50, 5).Value = '  -1.90%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.0505'
$ws.Cells.Item(51, 5).Value = '  -1.69%  '
